$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "64.469.65"
$ws.Range("E2").Value = "  -3.41%  "
$ws.Range("D3").Value = "3.398.86"
$ws.Range("E3").Value = "  -4.12%  "
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "578.30"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -4.83%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "131.98"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -9.09%  "
$ws.Range("E7").Value = "  -0.01%  "
$ws.Range("D8").Value = "3.401.36"
$ws.Range("E8").Value = "  -4.03%  "
$ws.Range("E9").Value = "  -7.42%  "
$ws.Range("E10").Value = "  -10.58%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.92"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -11.61%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.370"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -11.12%  "
$ws.Range("D13").Value = "3.975.57"
$ws.Range("E13").Value = "  -4.17%  "
$ws.Range("B14").Value = "ShibaInu"
$ws.Range("C14").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.0000175"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -11.17%  "
$ws.Range("B15").Value = "TRON"
$ws.Range("C15").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.115"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.85%  "
$ws.Range("B16").Value = "WrappedEther"
$ws.Range("C16").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D16").Value = "3.383.51"
$ws.Range("E16").Value = "  -4.99%  "
$ws.Range("D17").Value = "64.491.64"
$ws.Range("E17").Value = "  -3.23%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "25.80"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -11.57%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "9.31"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -15.74%  "
$ws.Range("E20").Value = "  -10.16%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.36"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -9.65%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "376.51"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -12.00%  "
$ws.Range("E23").Value = "  +0.07%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.537"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -10.66%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "71.35"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -8.26%  "
$ws.Range("D26").Value = "3.532.21"
$ws.Range("E26").Value = "  -4.24%  "
$ws.Range("E27").Value = "  -12.63%  "
$ws.Range("E28").Value = "  -0.15%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.07"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -12.17%  "
$ws.Range("E30").Value = "  -12.68%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.87"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -13.40%  "
$ws.Range("D32").Value = "3.416.64"
$ws.Range("E32").Value = "  -3.87%  "
$ws.Range("E34").Value = "  -7.10%  "
$ws.Range("E35").Value = "  -10.73%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "170.13"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -4.25%  "
$ws.Range("E37").Value = "  -14.23%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "6.57"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -14.65%  "
$ws.Range("E39").Value = "  -13.50%  "
$ws.Range("E40").Value = "  -14.65%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0752"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -9.59%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.792"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -8.49%  "
$ws.Range("E43").Value = "  -0.13%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "41.73"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -8.14%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "4.21"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -16.69%  "
$ws.Range("E46").Value = "  -12.17%  "
$ws.Range("E47").Value = "  -3.43%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "21.86"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -6.76%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "6.43"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -10.14%  "
$ws.Range("D50").Value = "2.178.46"
$ws.Range("E50").Value = "  -6.69%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "19.73"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -11.13%  "
